$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.895.26'
$ws.Range("E2").Value = '  -0.02%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.901.09'
$ws.Range("E3").Value = '  +0.67%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7632'
$ws.Range("E5").Value = '  +4.10%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '240.52'
$ws.Range("E6").Value = '  -0.75%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.001'
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3065'
$ws.Range("E8").Value = '  -1.29%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '25.53'
$ws.Range("E9").Value = '  -2.82%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06840'
$ws.Range("E10").Value = '  -0.90%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07978'
$ws.Range("E11").Value = '  +0.47%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.898.95'
$ws.Range("E12").Value = '  +0.32%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.7431'
$ws.Range("E13").Value = '  -3.00%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.159'
$ws.Range("E14").Value = '  -1.53%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '91.10'
$ws.Range("E15").Value = '  -0.17%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '29.889.51'
$ws.Range("E16").Value = '  -0.13%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '13.92'
$ws.Range("E17").Value = '  -1.66%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.934'
$ws.Range("E18").Value = '  +3.28%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '242.64'
$ws.Range("E19").Value = '  +1.18%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007697'
$ws.Range("E20").Value = '  -0.79%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.000'
$ws.Range("E21").Value = '  -0.10%  '
$ws.Range("E22").Value = '  +0.09%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.935'
$ws.Range("E23").Value = '  +0.32%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '166.77'
$ws.Range("E24").Value = '  +1.40%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.238'
$ws.Range("E25").Value = '  -0.74%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '18.70'
$ws.Range("E26").Value = '  -1.02%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1293'
$ws.Range("E27").Value = '  +1.80%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.033'
$ws.Range("E28").Value = '  +0.86%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.403'
$ws.Range("E29").Value = '  +3.48%  '
$ws.Range("E30").Value = '  -0.83%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.256'
$ws.Range("E31").Value = '  -1.10%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.075'
$ws.Range("E32").Value = '  -0.17%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05279'
$ws.Range("E33").Value = '  +3.79%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.252'
$ws.Range("E34").Value = '  -1.92%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7261'
$ws.Range("E35").Value = '  -1.43%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.715'
$ws.Range("E36").Value = '  -0.22%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01924'
$ws.Range("E37").Value = '  +0.03%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.785'
$ws.Range("E38").Value = '  +0.38%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.175'
$ws.Range("E39").Value = '  -2.44%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.4413'
$ws.Range("E40").Value = '  -0.84%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '71.98'
$ws.Range("E41").Value = '  -3.61%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.001'
$ws.Range("E42").Value = '  -0.03%  '
$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8311'
$ws.Range("E43").Value = '  -0.75%  '
$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.882'
$ws.Range("E44").Value = '  -2.48%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '7.620'
$ws.Range("E45").Value = '  +0.21%  '
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '9.790'
$ws.Range("E46").Value = '  +0.20%  '
$ws.Range("B47").Value = 'Quant'
$ws.Range("C47").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '100.05'
$ws.Range("E47").Value = '  -0.84%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.048.85'
$ws.Range("E48").Value = '  +0.73%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '36.10'
$ws.Range("E49").Value = '  -2.90%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.477'
$ws.Range("E50").Value = '  +1.58%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05942'
$ws.Range("E51").Value = '  +0.05%  '
